# Generate Report for handoff
# Updates the zh-cn and de-de status sheets: marks the file as ready for
# handoff, records the newly produced handoff (.xlf) file + timestamp, and
# switches the handoff reason from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$sheetsInfo = @(
  @{ Name = "zh-cn"; File = "1b32f26a-0aeb-43f4-9714-47b990a0f7d9.bb23dac6ce163c6ce7c5b26ccef762df990f6679.zh-cn.xlf"; Time = "2016-01-15 08:34:05" },
  @{ Name = "de-de"; File = "1b32f26a-0aeb-43f4-9714-47b990a0f7d9.bb23dac6ce163c6ce7c5b26ccef762df990f6679.de-de.xlf"; Time = "2016-01-15 08:34:14" }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Status: "Handoff transform failed" -> "Ready for handoff"
    $ws.Range("B2").Value = "Ready for handoff"

    # Latest Handoff File: new hyperlink cell pointing at the generated xlf
    $ws.Range("C2").Value = $info.File
    $ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/eaab7d9106d9b7b2229b1655240a4b306c3bcab9/e2e/" + $info.File, "", "", $info.File)
    $ws.Range("C2").Style = "HyperLink"

    # Latest Handoff Datetime
    $ws.Range("D2").Value = $info.Time

    # Handoff Reason: "Ignored" -> "Include"
    $ws.Range("H2").Value = "Include"
}
